# Natmi following Dr Hou advice
# Efna5-Ephb2 LR-pair sheet: add "ECs" as a possible Sending/Target cluster
# and refresh all computed statistic columns (E:T) for the full
# Sending-cluster x Target-cluster grid (FAPs,sCs) x (ECs,FAPs,sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data rows (2-7), columns A-T, in sheet order.
$rows = @(
    @("FAPs","Efna5","Ephb2","ECs",3,1,1.666083666666667,4.998251,0.6125276070882968,0.6125276070882968,1,0.3333333333333333,0.005966,0.017898,0.00125947234315407,0.00125947234315407,0.009939855155333332,0.089458696398,0.0007714615805460526,0.0007714615805460526),
    @("FAPs","Efna5","Ephb2","FAPs",3,1,1.666083666666667,4.998251,0.6125276070882968,0.6125276070882968,3,1,3.821776,11.465328,0.8068087787009701,0.8068087787009701,6.367398571258666,57.306587141328,0.4941926505955364,0.4941926505955364),
    @("FAPs","Efna5","Ephb2","sCs",3,1,1.666083666666667,4.998251,0.6125276070882968,0.6125276070882968,3,1,0.9091623333333333,2.727487,0.1919317489558758,0.1919317489558758,1.514740513915222,13.632664625237,0.1175634949122143,0.1175634949122143),
    @("sCs","Efna5","Ephb2","ECs",3,1,1.053930333333333,3.161791,0.3874723929117032,0.3874723929117031,1,0.3333333333333333,0.005966,0.017898,0.00125947234315407,0.00125947234315407,0.006287748368666667,0.056589735318,0.0004880107626080172,0.0004880107626080171),
    @("sCs","Efna5","Ephb2","FAPs",3,1,1.053930333333333,3.161791,0.3874723929117032,0.3874723929117031,3,1,3.821776,11.465328,0.8068087787009701,0.8068087787009701,4.027885653605334,36.250970882448,0.3126161281054337,0.3126161281054336),
    @("sCs","Efna5","Ephb2","sCs",3,1,1.053930333333333,3.161791,0.3874723929117032,0.3874723929117031,3,1,0.9091623333333333,2.727487,0.1919317489558758,0.1919317489558758,0.9581937610241111,8.623743849217,0.07436825404366147,0.07436825404366144)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
